$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$newValues = @(
    "48+27=",
    "16+56=",
    "47+18=",
    "58-39=",
    "31-13=",
    "40-9=",
    "78-9=",
    "92-48=",
    "28+7=",
    "34+49=",
    "92-33=",
    "17+64=",
    "58+5=",
    "68-39=",
    "72-68=",
    "70-63=",
    "92-16=",
    "7+48=",
    "75-39=",
    "39+4=",
    "80-63=",
    "60-44=",
    "14+49=",
    "28+4=",
    "19+7=",
    "54-8=",
    "46+15=",
    "85-78=",
    "39+18=",
    "42-34=",
    "35+28=",
    "81-39=",
    "18+74=",
    "92-3=",
    "21-19=",
    "90-69=",
    "25+36=",
    "54-25=",
    "73-58=",
    "84+7=",
    "8+48=",
    "7+48=",
    "33-17=",
    "72-56=",
    "90-87=",
    "8+49=",
    "25-19=",
    "44+7=",
    "85-68=",
    "82-79=",
    "27+37=",
    "92-29=",
    "82-78=",
    "90-45=",
    "82-79=",
    "49+45=",
    "53-24=",
    "37+56=",
    "4+8=",
    "18+26=",
    "5+57=",
    "92-9=",
    "22-17=",
    "40-27=",
    "87-39=",
    "27+67=",
    "18+8=",
    "6+59=",
    "38+15=",
    "89+6=",
    "7+65=",
    "46-28=",
    "73-5=",
    "48+14=",
    "53-25=",
    "68-9=",
    "44-5=",
    "37+29=",
    "71-42=",
    "29+54=",
    "15+17=",
    "26+35=",
    "98-49=",
    "19+37=",
    "50-36=",
    "15+59=",
    "61-28=",
    "51-12=",
    "61-15=",
    "37+26=",
    "52-45=",
    "39+52=",
    "49+13=",
    "51-3=",
    "96-78=",
    "78+3=",
    "35+6=",
    "14+57=",
    "94-69=",
    "38+58="
)

$numCols = $table.Columns.Count
$index = 0
foreach ($row in $table.Rows) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $row.Cells.Item($c)
        $cellRange = $cell.Range
        $moveResult = $cellRange.MoveEnd(1, -1)  # exclude the end-of-cell marker
        $cellRange.Text = $newValues[$index]
        $index = $index + 1
    }
}

Write-Output "Updated $index cells"
